# Generate Report for Handoff
# Updates the localization-status report to mark the
# fa92b67d-9c1f-47dc-9efa-f48ba314ec48.md file as handed off for
# translation (zh-cn and de-de), bumping its status from
# "In Translation" to "Ready for handoff", its priority from "ht" to
# "mt", and refreshing the relevant handoff timestamps.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

# --- Overview sheet: row 3 is the fa92b67d-... file ---
$overview.Range("E3").Value = "Ready for handoff"
$overview.Range("F3").Value = "Ready for handoff"
$overview.Range("G3").Value = "2016-08-16 20:11:37"

# --- zh-cn sheet: row 3 is the fa92b67d-... file ---
$zhcn.Range("C3").Value = "Ready for handoff"
$zhcn.Range("E3").Value = "mt"
$zhcn.Range("H3").Value = "2016-08-16 20:11:32"

# --- de-de sheet: row 3 is the fa92b67d-... file ---
$dede.Range("C3").Value = "Ready for handoff"
$dede.Range("E3").Value = "mt"
$dede.Range("H3").Value = "2016-08-16 20:11:37"

# Re-fit the columns that now hold the longer "Ready for handoff" text,
# matching Excel's autofit behaviour triggered by the new content.
# (16.333333333333336 is the COM ColumnWidth input that this host's
# pixel-snapping reproduces as ~17.22 characters wide.)
$overview.Range("E:F").ColumnWidth = 16.333333333333336
$zhcn.Range("C:C").ColumnWidth = 16.333333333333336
$dede.Range("C:C").ColumnWidth = 16.333333333333336
